$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102 updates
$ws.Range("B102").Value = 110.5
$ws.Range("C102").Value = 108.9
$ws.Range("D102").Value = 95.09999999999999
$ws.Range("E102").Value = 106.7
$ws.Range("F102").Value = 124.1
$ws.Range("G102").Value = 133
$ws.Range("H102").Value = 106
$ws.Range("I102").Value = 109.9
$ws.Range("J102").Value = 112.5

# Row 103 updates
$ws.Range("B103").Value = 115.1
$ws.Range("C103").Value = 102.7
$ws.Range("D103").Value = 99.7
$ws.Range("E103").Value = 108.2
$ws.Range("F103").Value = 101.2
$ws.Range("G103").Value = 130
$ws.Range("H103").Value = 120.2
$ws.Range("I103").Value = 114.3
$ws.Range("J103").Value = 117.1
